$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet2 ("A1:O98") - add/update the fitted "B" column values (percentages)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Rows whose B value changed (already had a value/style before)
$changedRows = @{
    2 = 0.3
    3 = 0.33
    4 = 0.4
    5 = 0.43
    7 = 0.5
    8 = 0.57999999999999996
    9 = 0.6
}

# Rows whose B cell is brand new (needs value + percentage number format)
$newRows = @{
    14 = 0.7
    24 = 0.84
    25 = 0.84
    26 = 0.84
    27 = 0.84
    28 = 0.84
    30 = 0.85
    31 = 0.85
    32 = 0.85
    33 = 0.85
    34 = 0.85
    35 = 0.85
    36 = 0.85
    37 = 0.85
    38 = 0.85
    40 = 0.86
    41 = 0.87
    43 = 0.88
    44 = 0.88
    45 = 0.89
    46 = 0.89
    48 = 0.91
    49 = 0.92
    50 = 0.93
    51 = 0.94
    52 = 0.94
    54 = 0.95
    55 = 0.95
    56 = 0.95
    57 = 0.95
    58 = 0.95
    60 = 0.96
    61 = 0.96
    62 = 0.96
    63 = 0.96
    64 = 0.96
    66 = 0.98
    67 = 0.98
    68 = 0.98
    69 = 0.98
    70 = 0.98
    71 = 0.98
    72 = 0.98
    73 = 0.98
    74 = 0.98
    75 = 0.98
    76 = 0.98
    77 = 0.98
    78 = 0.98
    79 = 0.98
    80 = 0.98
    81 = 0.98
    82 = 0.98
    83 = 0.98
    84 = 0.98
    85 = 0.98
    86 = 0.98
    87 = 0.98
    88 = 0.98
    89 = 0.98
    90 = 0.98
    91 = 0.99
    92 = 0.99
    93 = 0.99
    94 = 0.99
    95 = 0.99
    96 = 0.99
    97 = 0.99
}

foreach ($row in $changedRows.Keys) {
    $ws2.Cells.Item($row, 2).Value = $changedRows[$row]
}

foreach ($row in $newRows.Keys) {
    $ws2.Cells.Item($row, 2).Value = $newRows[$row]
    $ws2.Cells.Item($row, 2).NumberFormat = "0%"
}

# Make Sheet2 the active sheet/tab and select D3 on it (matches the new
# tabSelected + selection stored in the workbook).
$ws2.Activate() | Out-Null
$ws2.Range("D3").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet3 ("A1:U9") - page setup now specifies paper size / orientation
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1
